$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 59 held phone "09876543" (stored as text, with its leading
# zero preserved) together with 0 points. The update keeps that original
# text record intact but pushes it down to row 60, and inserts a new row
# 59 above it holding the numeric, zero-padding-stripped phone value
# 9876543 with 0 points (points 09876543 -> 0.00).

# Insert a blank row above row 59; this shifts the old row 59 down to 60.
$ws.Rows.Item(59).Insert()

# Row 59's birthday column (B) stays blank, just like the row that was
# pushed down to 60. Copy that empty cell down so row 59 gets a real
# (but empty) cell there instead of being left out entirely.
$ws.Cells.Item(60, 2).Copy($ws.Cells.Item(59, 2))

# Fill in the new row 59: numeric phone 9876543 and 0 points.
$ws.Cells.Item(59, 1).Value = 9876543
$ws.Cells.Item(59, 3).Value = 0
